# Add files via upload
# - Updates row 10/11 data + averages on "Proportional - Proportional"
# - Adds a new "Proportional - Random" sheet (after the last sheet) with its
#   own Test/Nodes/Edges/Incentive data + averages, and makes it the active sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Fix up the last two data rows (and the dependent AVERAGE formulas) on the
#    "Proportional - Proportional" sheet.
# ---------------------------------------------------------------------------
$wsPP = $wb.Worksheets.Item("Proportional - Proportional")

$wsPP.Range("C10").Value = 59881
$wsPP.Range("D10").Value = 36349
$wsPP.Range("E10").Value = 134704

$wsPP.Range("C11").Value = 59869
$wsPP.Range("D11").Value = 36358
$wsPP.Range("E11").Value = 134716

$wsPP.Range("D12").Formula = "=AVERAGE(D2:D11)"
$wsPP.Range("E12").Formula = "=AVERAGE(E2:E11)"

# that sheet is no longer the selected tab / its selection moves to E12
[void]$wsPP.Range("E12").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "Proportional - Random" worksheet as the last tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPR = $wb.Worksheets.Add($null, $lastSheet)
$wsPR.Name = "Proportional - Random"

$wsPR.Range("A1").Value = "Test"
$wsPR.Range("B1").Value = "Nodes"
$wsPR.Range("C1").Value = "Edges"
$wsPR.Range("D1").Value = "Number of Nodes with Incentive"
$wsPR.Range("E1").Value = "Total Incentive"

$rows = @(
    @(1,  37873, 59874, 27007, 65536),
    @(2,  37873, 59891, 27008, 65135),
    @(3,  37873, 59893, 26900, 65230),
    @(4,  37873, 59866,     0,     0),
    @(5,  37873,     0,     0,     0),
    @(6,  37873,     0,     0,     0),
    @(7,  37873,     0,     0,     0),
    @(8,  37873,     0,     0,     0),
    @(9,  37873,     0,     0,     0),
    @(10, 37873,     0,     0,     0)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $wsPR.Cells.Item($r, 1).Value = $rows[$i][0]
    $wsPR.Cells.Item($r, 2).Value = $rows[$i][1]
    $wsPR.Cells.Item($r, 3).Value = $rows[$i][2]
    $wsPR.Cells.Item($r, 4).Value = $rows[$i][3]
    $wsPR.Cells.Item($r, 5).Value = $rows[$i][4]
}

$wsPR.Range("C12").Value = "Media"
$wsPR.Range("D12").Formula = "=AVERAGE(D2:D11)"
$wsPR.Range("E12").Formula = "=AVERAGE(E2:E11)"

# column widths matching the other sheets in the workbook
$wsPR.Columns.Item(4).ColumnWidth = 27.77734375
$wsPR.Columns.Item(5).ColumnWidth = 13.21875

[void]$wsPR.Range("D5").Select()
[void]$wsPR.Activate()
